$wb = $excel.ActiveWorkbook

# Update the Selex sheet values for the sex structure run
$selex = $wb.Worksheets.Item("Selex")
$selex.Range("B2").Value = 0.6
$selex.Range("B3").Value = 62.5
$selex.Range("B4").Value = 0.8
$selex.Range("B5").Value = 52.5

# Update the selection on the Controls sheet
$controls = $wb.Worksheets.Item("Controls")
$controls.Activate()
$controls.Range("B3").Select()
